$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2.655341483754701
$ws.Range("E2").Value = 11.261487360959372
$ws.Range("F2").Value = 0.10698582780178503

# Row 3
$ws.Range("D3").Value = 3.7766184336060604
$ws.Range("E3").Value = 17.07731098594899
$ws.Range("F3").Value = 0.27960155269241754

# Row 4
$ws.Range("D4").Value = 2.415386042556059
$ws.Range("E4").Value = 12.682530365509892
$ws.Range("F4").Value = 0.2838364403663997

# Row 5
$ws.Range("D5").Value = 2.40895930679867
$ws.Range("E5").Value = 10.97313951670336
$ws.Range("F5").Value = 0.08923200269448178

# Row 6
$ws.Range("D6").Value = 0.9497513267881987
$ws.Range("E6").Value = 4.5330957394836995
$ws.Range("F6").Value = 0.04336517260314411
$ws.Range("G6").Value = 0.0006
$ws.Range("H6").Value = 0.0036

# Row 7
$ws.Range("D7").Value = 2.6326877789750918
$ws.Range("E7").Value = 15.454037717279641
$ws.Range("F7").Value = 0.23610518550484821
